# Rename the header columns from the old "_old"/"_new" suffix convention
# to the new format-version based suffixes ("_FV2404" / "_FV2410").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("_old", "_FV2404", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$ws.Cells.Replace("_new", "_FV2410", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

# Turn the data range into a real Excel Table (ListObject) so it gets a
# header row / autofilter / structured reference support.
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
